# Auto-generated PowerShell Excel COM-interop script
# Applies the "may 9th" data change: inserts 9 new rows of sensor readings
# right after the header row (shifting the original data down), and appends
# one more new row of sensor readings at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 blank rows right after the header (before the existing data),
# pushing the existing data rows down from row 2-21 to row 11-30.
$ws.Rows("2:10").Insert()

# The inserted rows pick up formatting from the row above (the header);
# clear that so the new data rows are unstyled, like the rest of the table.
$ws.Range("A2:H10").ClearFormats()

# Fill the newly inserted rows with the new sensor data (timestamps 0-800).
$newTopRows = @(
    @(0, -0.5562429428100586, 1.515548229217529, 0.1288182139396667, -0.007375299738829639, -0.04188020327402392, 0.02404832670136416),
    @(100, -0.7228193283081055, 1.569920063018799, 0.3090478777885437, 0.009710959871025781, -0.04103577224647274, -0.004662338863401103),
    @(200, -0.7100648880004883, 1.677208662033081, 0.4076560139656067, 0.04170951860792497, 0.04793495536946196, -0.07793023174299907),
    @(300, -0.8624534606933594, 1.662384271621704, 0.1703254878520965, -0.02235946409842538, -0.01261257046066655, 0.04359601744834121),
    @(400, -0.2973442077636719, 1.564767122268677, 0.0300358235836029, -0.03681361302733401, -0.006297301829737623, -0.01943090470398167),
    @(500, -0.7404184341430664, 1.651389360427856, 0.2246546447277069, 0.02995036389021311, 0.01477754981640481, -0.02054483487325551),
    @(600, -0.7922754287719727, 1.612479209899902, 0.1184005141258239, 0.04258089907029088, -0.09914881779867053, -0.03582545123336939),
    @(700, -0.5356760025024414, 1.573039531707764, 0.1370119750499725, 0.0959058403968811, -0.1950187236070633, 0.0042760567739605),
    @(800, -0.6788949966430664, 1.517318725585938, 0.1871603727340698, -0.04207783586838668, -0.3005187625394148, 0.05587620359352399),
)

$r = 2
foreach ($row in $newTopRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = "falling"
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# Append a brand-new row (timestamp 2900) at the end of the table.
$ws.Range("A31").Value = 2900
$ws.Range("B31").Value = "falling"
$ws.Range("C31").Value = -0.1529102325439453
$ws.Range("D31").Value = 1.506775379180908
$ws.Range("E31").Value = 0.2578078508377075
$ws.Range("F31").Value = -0.009441461181268051
$ws.Range("G31").Value = 0.03861925794797785
$ws.Range("H31").Value = -0.04366788180435401
